$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "guia" numbers that replace the old A2:A40 block (corrected data load).
$newValues = @(
    2129031526,2129034230,2129036584,2139595988,2145085871,2151145588,2152483425,
    2158670584,2164814489,2167929154,2167929401,2171535321,2171535509,2171921928,
    2173240654,2173785485,2176081466,2176738241,2176746209,2176746924,2176747539,
    2176751705,2176779765,2176901371,2176938657,2176938663,2179444401,3011003965,
    3014466745,3020072394,3041022504,3043847040,3045539096,3047099382,9156674647,
    9159055187,9159793527,9159932900,9161309827
)

# Drop the old rows 41:70 entirely (shifting rows up) so the stale values
# (including the old A41:A55 block that no longer exists afterwards) are gone.
$ws.Rows("41:70").Delete()

# Write the corrected values into A2:A40.
for ($i = 0; $i -lt $newValues.Length; $i++) {
    $ws.Cells.Item($i + 2, 1).Value = $newValues[$i]
}

# Recreate the trailing placeholder rows 56:70 as empty cells that keep the
# original right-aligned style (s="1", copied from a still-styled A2:A40
# cell) without resurrecting rows 41:55.
$ws.Range("A2").Copy()
$ws.Range("A56:A70").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Strip the old right-aligned style (s="1") from A2:A40 so they fall back to
# the default cell style, matching the refreshed data dump.
$ws.Range("A2:A40").ClearFormats()

# Match the saved selection state (A2 active, A2:A40 selected).
$ws.Range("A2:A40").Select()
